# Applies updated market-price / profit figures to the Leve profit tracking sheets.
# Generated to match the upstream data refresh (scheduled Sheets runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 4 - Root Rush
$ws.Cells.Item(4, 8).Value = 87.5
$ws.Cells.Item(4, 9).Value = 87.5
$ws.Cells.Item(4, 11).Value = 87.5
$ws.Cells.Item(4, 13).Value = 26.5
# row 64 - Forged from the Void
$ws.Cells.Item(64, 8).Value = 7542.346
$ws.Cells.Item(64, 9).Value = 6304.6
$ws.Cells.Item(64, 11).Value = 6304.6
$ws.Cells.Item(64, 13).Value = -6056.6
# row 67 - Dodging the Draft (L)
$ws.Cells.Item(67, 8).Value = 7542.346
$ws.Cells.Item(67, 9).Value = 6304.6
$ws.Cells.Item(67, 11).Value = 6304.6
$ws.Cells.Item(67, 13).Value = -5446.6
# row 98 - The Dotted Line
$ws.Cells.Item(98, 8).Value = 1590
$ws.Cells.Item(98, 9).Value = 1487.5
$ws.Cells.Item(98, 11).Value = 1487.5
$ws.Cells.Item(98, 13).Value = 10.5
# row 100 - Asking for a Friend
$ws.Cells.Item(100, 8).Value = 5670.65
$ws.Cells.Item(100, 9).Value = 3344.375
$ws.Cells.Item(100, 11).Value = 3344.375
$ws.Cells.Item(100, 13).Value = -2803.375
# row 113 - Amaro Kart
$ws.Cells.Item(113, 8).Value = 5427.7417
$ws.Cells.Item(113, 9).Value = 4152.467
$ws.Cells.Item(113, 10).Value = 6623.3125
$ws.Cells.Item(113, 11).Value = 4152.467
$ws.Cells.Item(113, 12).Value = 6623.3125
$ws.Cells.Item(113, 13).Value = -898.4669999999996
$ws.Cells.Item(113, 14).Value = -13131.3125
# row 122 - Wishful Inking
$ws.Cells.Item(122, 8).Value = 1590
$ws.Cells.Item(122, 9).Value = 1487.5
$ws.Cells.Item(122, 11).Value = 4462.5
$ws.Cells.Item(122, 13).Value = -2012.5
# row 137 - Cutting Edge of Culinary Quality
$ws.Cells.Item(137, 8).Value = 2985
$ws.Cells.Item(137, 9).Value = 3646.6667
$ws.Cells.Item(137, 11).Value = 10940.0001
$ws.Cells.Item(137, 13).Value = -8390.000100000001

$ws = $wb.Worksheets.Item("ARM")
# row 88 - The Mast Chance
$ws.Cells.Item(88, 8).Value = 1770.125
$ws.Cells.Item(88, 9).Value = 1845.1111
$ws.Cells.Item(88, 10).Value = 1673.7142
$ws.Cells.Item(88, 11).Value = 1845.1111
$ws.Cells.Item(88, 12).Value = 1673.7142
$ws.Cells.Item(88, 13).Value = -1439.1111
$ws.Cells.Item(88, 14).Value = -2485.7142
# row 91 - The Rose and the Riveter (L)
$ws.Cells.Item(91, 8).Value = 1770.125
$ws.Cells.Item(91, 9).Value = 1845.1111
$ws.Cells.Item(91, 10).Value = 1673.7142
$ws.Cells.Item(91, 11).Value = 1845.1111
$ws.Cells.Item(91, 12).Value = 1673.7142
$ws.Cells.Item(91, 13).Value = -441.1111000000001
$ws.Cells.Item(91, 14).Value = -4481.7142
# row 97 - Ore for Me
$ws.Cells.Item(97, 8).Value = 1174.6666
$ws.Cells.Item(97, 9).Value = 1347.5834
$ws.Cells.Item(97, 11).Value = 1347.5834
$ws.Cells.Item(97, 13).Value = -851.5834
# row 110 - Scheduled Maintenance
$ws.Cells.Item(110, 8).Value = 2611.4666
$ws.Cells.Item(110, 9).Value = 2611.4666
$ws.Cells.Item(110, 11).Value = 2611.4666
$ws.Cells.Item(110, 13).Value = -566.4666000000002

$ws = $wb.Worksheets.Item("BSM")
# row 22 - Riveting Run
$ws.Cells.Item(22, 8).Value = 4376.2
$ws.Cells.Item(22, 9).Value = 4376.2
$ws.Cells.Item(22, 11).Value = 4376.2
$ws.Cells.Item(22, 13).Value = -4203.2
# row 26 - Unseamly Conditions
$ws.Cells.Item(26, 8).Value = 0
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 13).Value = $null
# row 74 - I Could Feel That from Here
$ws.Cells.Item(74, 8).Value = 54999.75
$ws.Cells.Item(74, 10).Value = 70999.5
$ws.Cells.Item(74, 12).Value = 70999.5
$ws.Cells.Item(74, 14).Value = -72871.5
# row 77 - Tensions in Creasing (L)
$ws.Cells.Item(77, 8).Value = 54999.75
$ws.Cells.Item(77, 10).Value = 70999.5
$ws.Cells.Item(77, 12).Value = 212998.5
$ws.Cells.Item(77, 14).Value = -222358.5
# row 105 - Ingot to Wing It
$ws.Cells.Item(105, 8).Value = 3844.2104
$ws.Cells.Item(105, 9).Value = 2467.3635
$ws.Cells.Item(105, 11).Value = 2467.3635
$ws.Cells.Item(105, 13).Value = -720.3634999999999

$ws = $wb.Worksheets.Item("CRP")
# row 16 - Raise the Roof
$ws.Cells.Item(16, 8).Value = 519.6429000000001
$ws.Cells.Item(16, 9).Value = 519.6429000000001
$ws.Cells.Item(16, 11).Value = 519.6429000000001
$ws.Cells.Item(16, 13).Value = -232.6429000000001
# row 107 - Built to Last
$ws.Cells.Item(107, 8).Value = 3334.4443
$ws.Cells.Item(107, 9).Value = 458.86957
$ws.Cells.Item(107, 10).Value = 8422
$ws.Cells.Item(107, 11).Value = 458.86957
$ws.Cells.Item(107, 12).Value = 8422
$ws.Cells.Item(107, 13).Value = 1461.13043
$ws.Cells.Item(107, 14).Value = -12262
# row 109 - Playing the Market
$ws.Cells.Item(109, 8).Value = 41056.5
$ws.Cells.Item(109, 10).Value = 41056.5
$ws.Cells.Item(109, 12).Value = 41056.5
$ws.Cells.Item(109, 14).Value = -43136.5
# row 113 - Patient Patients
$ws.Cells.Item(113, 8).Value = 519.6429000000001
$ws.Cells.Item(113, 9).Value = 519.6429000000001
$ws.Cells.Item(113, 11).Value = 519.6429000000001
$ws.Cells.Item(113, 13).Value = 1650.3571

$ws = $wb.Worksheets.Item("CUL")
# row 14 - Keep Your Powder Dry
$ws.Cells.Item(14, 8).Value = 176.46666
$ws.Cells.Item(14, 9).Value = 176.46666
$ws.Cells.Item(14, 11).Value = 529.3999799999999
$ws.Cells.Item(14, 13).Value = -356.3999799999999
# row 107 - Slippery Service
$ws.Cells.Item(107, 8).Value = 305.15384
$ws.Cells.Item(107, 9).Value = 391.66666
$ws.Cells.Item(107, 11).Value = 1174.99998
$ws.Cells.Item(107, 13).Value = 745.0000199999999
# row 113 - Can't Eat Just One
$ws.Cells.Item(113, 8).Value = 1215.037
$ws.Cells.Item(113, 9).Value = 984.0909
$ws.Cells.Item(113, 11).Value = 2952.2727
$ws.Cells.Item(113, 13).Value = -782.2727

$ws = $wb.Worksheets.Item("GSM")
# row 99 - Needle in a Hingan Stack
$ws.Cells.Item(99, 8).Value = 9686.5
$ws.Cells.Item(99, 9).Value = 9686.5
$ws.Cells.Item(99, 11).Value = 9686.5
$ws.Cells.Item(99, 13).Value = -7440.5
# row 132 - On Board for Lar
$ws.Cells.Item(132, 8).Value = 4485.2256
$ws.Cells.Item(132, 10).Value = 4247.846
$ws.Cells.Item(132, 12).Value = 12743.538
$ws.Cells.Item(132, 14).Value = -17803.538

$ws = $wb.Worksheets.Item("LTW")
# row 7 - Tan Before the Ban
$ws.Cells.Item(7, 8).Value = 5458.6924
$ws.Cells.Item(7, 9).Value = 3411.6
$ws.Cells.Item(7, 10).Value = 6738.125
$ws.Cells.Item(7, 11).Value = 3411.6
$ws.Cells.Item(7, 12).Value = 6738.125
$ws.Cells.Item(7, 13).Value = -3299.6
$ws.Cells.Item(7, 14).Value = -6962.125
# row 39 - Quality over Quantity
$ws.Cells.Item(39, 8).Value = 4752
$ws.Cells.Item(39, 9).Value = 4752
$ws.Cells.Item(39, 11).Value = 4752
$ws.Cells.Item(39, 13).Value = -4292
# row 93 - Hide to Go Seek
$ws.Cells.Item(93, 8).Value = 272263.4
$ws.Cells.Item(93, 9).Value = 2001
$ws.Cells.Item(93, 11).Value = 2001
$ws.Cells.Item(93, 13).Value = -753
# row 100 - Tiger in the Sack
$ws.Cells.Item(100, 8).Value = 53684.453
$ws.Cells.Item(100, 9).Value = 75737.2
$ws.Cells.Item(100, 11).Value = 75737.2
$ws.Cells.Item(100, 13).Value = -75196.2
# row 122 - Hell on Leather
$ws.Cells.Item(122, 8).Value = 5747.737
$ws.Cells.Item(122, 9).Value = 4486
$ws.Cells.Item(122, 11).Value = 13458
$ws.Cells.Item(122, 13).Value = -11008
# row 126 - Battered Books
$ws.Cells.Item(126, 8).Value = 5458.6924
$ws.Cells.Item(126, 9).Value = 3411.6
$ws.Cells.Item(126, 10).Value = 6738.125
$ws.Cells.Item(126, 11).Value = 10234.8
$ws.Cells.Item(126, 12).Value = 20214.375
$ws.Cells.Item(126, 13).Value = -7764.799999999999
$ws.Cells.Item(126, 14).Value = -25154.375
# row 135 - Dreams of Ja
$ws.Cells.Item(135, 8).Value = 82725.30499999999
$ws.Cells.Item(135, 10).Value = 82725.30499999999
$ws.Cells.Item(135, 12).Value = 82725.30499999999
$ws.Cells.Item(135, 14).Value = -92865.30499999999
# row 136 - Respect for Br'aax
$ws.Cells.Item(136, 8).Value = 4358.8223
$ws.Cells.Item(136, 9).Value = 2332.6667
$ws.Cells.Item(136, 10).Value = 6674.4287
$ws.Cells.Item(136, 11).Value = 6998.000100000001
$ws.Cells.Item(136, 12).Value = 20023.2861
$ws.Cells.Item(136, 13).Value = -4448.000100000001
$ws.Cells.Item(136, 14).Value = -25123.2861

$ws = $wb.Worksheets.Item("WVR")
# row 99 - Say Yes to Formal Dress
$ws.Cells.Item(99, 8).Value = 32216
$ws.Cells.Item(99, 9).Value = 30432
$ws.Cells.Item(99, 11).Value = 30432
$ws.Cells.Item(99, 13).Value = -27437
# row 107 - Flax Wax
$ws.Cells.Item(107, 8).Value = 618.4545000000001
$ws.Cells.Item(107, 9).Value = 530.8333
$ws.Cells.Item(107, 10).Value = 723.6
$ws.Cells.Item(107, 11).Value = 1592.4999
$ws.Cells.Item(107, 12).Value = 2170.8
$ws.Cells.Item(107, 13).Value = 327.5001
$ws.Cells.Item(107, 14).Value = -6010.8
# row 113 - A Tender Table
$ws.Cells.Item(113, 8).Value = 615.2059
$ws.Cells.Item(113, 9).Value = 378.8846
$ws.Cells.Item(113, 11).Value = 1136.6538
$ws.Cells.Item(113, 13).Value = 1033.3462
# row 122 - Heavy Armoire
$ws.Cells.Item(122, 8).Value = 5809.7393
$ws.Cells.Item(122, 9).Value = 2701.0557
$ws.Cells.Item(122, 11).Value = 8103.1671
$ws.Cells.Item(122, 13).Value = -5653.1671
